$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on every cell we are about to write so that
# numeric-looking strings (e.g. "609.10") are stored as text, matching the
# inline/shared string cells in the source workbook rather than being
# auto-converted to floating point numbers by Excel.

$dCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D35", "D36", "D37", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '72.957.17'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '3.972.87'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '609.10'
$ws.Range("E5").Value = '  +4.69%  '
$ws.Range("D6").Value = '173.38'
$ws.Range("E6").Value = '  +13.74%  '
$ws.Range("D7").Value = '0.690'
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.798'
$ws.Range("E9").Value = '  +5.07%  '
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +8.51%  '
$ws.Range("D11").Value = '57.57'
$ws.Range("E11").Value = '  +6.78%  '
$ws.Range("D12").Value = '0.0000334'
$ws.Range("E12").Value = '  +2.58%  '
$ws.Range("D13").Value = '11.75'
$ws.Range("E13").Value = '  +5.65%  '
$ws.Range("D14").Value = '4.604.60'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").Value = '3.979.10'
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").Value = '21.32'
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '1.25'
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '14.22'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '72.938.08'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '453.67'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").Value = '4.87'
$ws.Range("E22").Value = '  +5.79%  '
$ws.Range("D23").Value = '97.30'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -4.40%  '
$ws.Range("D25").Value = '14.29'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").Value = '4.27'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").Value = '11.43'
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("D28").Value = '10.53'
$ws.Range("E28").Value = '  -4.07%  '
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").Value = '36.18'
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("E31").Value = '  +1.94%  '
$ws.Range("D32").Value = '14.03'
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").Value = '49.93'
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("D35").Value = '0.0000102'
$ws.Range("E35").Value = '  +15.09%  '
$ws.Range("D36").Value = '69.42'
$ws.Range("E36").Value = '  +2.90%  '
$ws.Range("D37").Value = '635.87'
$ws.Range("E37").Value = '  -7.68%  '
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").Value = '3.26'
$ws.Range("E43").Value = '  +46.41%  '
$ws.Range("D44").Value = '0.0487'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").Value = '10.57'
$ws.Range("E45").Value = '  -6.67%  '
$ws.Range("D46").Value = '0.150'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").Value = '2.99'
$ws.Range("E47").Value = '  -10.77%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").Value = '2.67'
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").Value = '0.000294'
$ws.Range("E49").Value = '  +7.98%  '
$ws.Range("D50").Value = '3.42'
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '149.62'
$ws.Range("E51").Value = '  +4.20%  '

# Restore default (no explicit number format) styling on the cells we
# touched so the saved XML matches the original "no s attribute" cells.
foreach ($ref in $dCells) {
    $ws.Range($ref).Style = "Normal"
}

Write-Output "Applied all cell updates"
